$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. delivery_man table: "status" field value-type/description updated
#    (row 15: C15 int8 -> varchar20, D15 "status" -> full status description)
# ---------------------------------------------------------------------------
$ws.Range("C15").Value2 = "varchar20"
$ws.Range("D15").Value2 = 'status "IDLE" or "Processing"'

# ---------------------------------------------------------------------------
# 2. order table: new "price" field added as the last row of that table
#    (a new row is inserted right before the old row 40, pushing every
#    following row down by one)
# ---------------------------------------------------------------------------
$ws.Rows.Item(40).Insert()

# Copy the formatting (styles/borders) from the row above so the new row
# matches the rest of the table instead of picking up a blank style.
$ws.Range("A39:D39").Copy()
$ws.Range("A40:D40").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

$ws.Range("B40").Value2 = "price"
$ws.Range("C40").Value2 = "double"
$ws.Range("D40").Value2 = "order price"

# The "order" table's merged label cell (A18:A39) now needs to cover the
# newly added row as well.
$null = $ws.Range("A18:A39").UnMerge()
$null = $ws.Range("A18:A40").Merge()

# Merging re-applies a top/middle/bottom bordered style triple to the range;
# restore the single uniform style the column used before (same one still
# used by every other untouched row in the table).
$ws.Range("A17").Copy()
$ws.Range("A18:A40").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Selection / scroll position left by the author after the edit
# ---------------------------------------------------------------------------
$null = $ws.Range("D41").Select()
